$d = $word.ActiveDocument

# Update the date in the first paragraph
$d.Content.Find.Execute("2025-12-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-27 Saturday", 2) | Out-Null

# Update the multiplication answers in the table, cell by cell (row, col)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "41×93=3813"
$t.Cell(1, 2).Range.Text = "79×46=3634"
$t.Cell(1, 3).Range.Text = "16×95=1520"
$t.Cell(1, 4).Range.Text = "73×11=803"
$t.Cell(1, 5).Range.Text = "74×85=6290"

$t.Cell(5, 1).Range.Text = "21×42=882"
$t.Cell(5, 2).Range.Text = "37×47=1739"
$t.Cell(5, 3).Range.Text = "55×32=1760"
$t.Cell(5, 4).Range.Text = "11×27=297"
$t.Cell(5, 5).Range.Text = "17×60=1020"

$t.Cell(10, 1).Range.Text = "93×86=7998"
$t.Cell(10, 2).Range.Text = "71×36=2556"
$t.Cell(10, 3).Range.Text = "59×12=708"
$t.Cell(10, 4).Range.Text = "68×11=748"
$t.Cell(10, 5).Range.Text = "38×45=1710"

$t.Cell(15, 1).Range.Text = "19×76=1444"
$t.Cell(15, 2).Range.Text = "17×32=544"
$t.Cell(15, 3).Range.Text = "27×44=1188"
$t.Cell(15, 4).Range.Text = "36×35=1260"
$t.Cell(15, 5).Range.Text = "33×75=2475"

$t.Cell(20, 1).Range.Text = "99×50=4950"
$t.Cell(20, 2).Range.Text = "17×76=1292"
$t.Cell(20, 3).Range.Text = "16×59=944"
$t.Cell(20, 4).Range.Text = "37×66=2442"
$t.Cell(20, 5).Range.Text = "99×87=8613"

